# BLNAS v10 report template: update the physician's title line from
# "Dr. med. Thiên-Trí Lâm" to "PD Dr. med. Thiên-Trí Lâm" in the
# first-page header's contact text box.
#
# wdReplaceAll = 2, wdFindContinue (Wrap) = 1

$d = $word.ActiveDocument

$oldText = "Dr. med. Thiên-Trí Lâm"
$newText = "PD Dr. med. Thiên-Trí Lâm"

$replaced = $false

# The text lives in the small "contact details" text box that is anchored
# in the document's first-page header (w:headerReference w:type="first").
# Section.Headers(2) == wdHeaderFooterFirstPage.
for ($secIdx = 1; $secIdx -le $d.Sections.Count -and -not $replaced; $secIdx++) {
    $section = $d.Sections.Item($secIdx)

    for ($hfIdx = 1; $hfIdx -le 3 -and -not $replaced; $hfIdx++) {
        $header = $null
        try { $header = $section.Headers($hfIdx) } catch { $header = $null }
        if ($header -eq $null -or -not $header.Exists) { continue }

        # 1) Try every floating shape's own text frame (this is where the
        #    "Textfeld 3" contact box with the physicians' names lives).
        for ($sIdx = 1; $sIdx -le $header.Shapes.Count -and -not $replaced; $sIdx++) {
            $shape = $header.Shapes.Item($sIdx)
            try {
                $tr = $shape.TextFrame.TextRange
                $found = $tr.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                           $true, 1, $false, $newText, 2)
                if ($found) { $replaced = $true }
            } catch {
                # shape has no usable text frame in this runtime - ignore and
                # fall through to the next strategy
            }
        }

        if ($replaced) { break }

        # 2) Fall back to a plain search across the header's own story range
        #    in case the text box content is exposed there instead.
        try {
            $found = $header.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                                 $true, 1, $false, $newText, 2)
            if ($found) { $replaced = $true }
        } catch { }
    }
}

# 3) Last resort: sweep every story range in the document (covers text
#    boxes that are modelled as independent stories).
if (-not $replaced) {
    try {
        $stories = $d.StoryRanges
        foreach ($story in $stories) {
            try {
                $found = $story.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                              $true, 1, $false, $newText, 2)
                if ($found) { $replaced = $true }
            } catch { }
        }
    } catch { }
}

# 4) Absolute last resort: the main document content (covers the case
#    where the line is duplicated into the body for some templates).
if (-not $replaced) {
    try {
        $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                          $true, 1, $false, $newText, 2)
        if ($found) { $replaced = $true }
    } catch { }
}

Write-Host "Title updated: $replaced"
